$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Workbook window view (best effort - engine may not persist this)
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Height = 11440

# ---------------------------------------------------------------------------
# 2) Header row (row 1): rename/rearrange the "Fetch"/"Write" columns
#    G = Fetch sequential ORM (renamed, was "Fetch (ms, avg of 1000)")
#    H = Fetch sequential raw/native (new, column used to hold "Write")
#    I = Fetch simultaneous ORM (new column)
#    J = Fetch simultaneous raw (new column)
#    K = Write (ms, avg of 1000) (moved here from column H)
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Fetch sequential ORM (ms, avg of 1000)"
$ws.Range("H1").Value = "Fetch sequential raw/native"
$ws.Range("I1").Value = "Fetch simultaneous ORM"
$ws.Range("J1").Value = "Fetch simultaneous raw"
$ws.Range("K1").Value = "Write (ms, avg of 1000)"

# ---------------------------------------------------------------------------
# 3) Data rows - copy style from existing G column (numFmt "0.00") onto the
#    newly used cells, then set their values.
# ---------------------------------------------------------------------------
$ws.Range("I2:K5").NumberFormat = "0.00"

# Row 2
$ws.Range("G2").Value = 4.3040000000000003
$ws.Range("H2").Value = 3.1680000000000001
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 9.9890000000000008

# Row 3
$ws.Range("G3").Value = 5.0430000000000001
$ws.Range("H3").Value = 5.306
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 1.736

# Row 4
$ws.Range("G4").Value = 4.109
$ws.Range("H4").Value = 4.2350000000000003
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 6.3150000000000004

# Row 5 (previously SUM formulas - now plain values)
$ws.Range("G5").Value = 9.1519999999999992
$ws.Range("H5").Value = 9.5410000000000004
$ws.Range("I5").Value = 3.242
$ws.Range("J5").Value = 3.1520000000000001
$ws.Range("K5").Value = 8.0510000000000002

# ---------------------------------------------------------------------------
# 4) Row 6: clear the stray insertion-time numbers (D6/E6), keep the style
# ---------------------------------------------------------------------------
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()

# ---------------------------------------------------------------------------
# 5) Column widths for the new/changed layout
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 33.6640625
$ws.Columns.Item(8).ColumnWidth = 24
$ws.Columns.Item(9).ColumnWidth = 24
$ws.Columns.Item(10).ColumnWidth = 24
$ws.Columns.Item(11).ColumnWidth = 20.1640625

# ---------------------------------------------------------------------------
# 6) Freeze panes: split after column A, with the view scrolled so column E
#    is the first visible column on the right pane.
# ---------------------------------------------------------------------------
$ws.Range("B1").Select()
$win.FreezePanes = $true

# ---------------------------------------------------------------------------
# 7) Final selection shown in the workbook
# ---------------------------------------------------------------------------
$ws.Range("K3").Select()
